# Add a "Save" column (H) to the s_vals sheet, matching the header style
# of the existing columns (B1:G1) and filling H2:H4 with 1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell
$ws.Range("H1").Value = "Save"

# Copy the formatting (bold, border, centered) from the neighboring
# header cell so the new header matches the existing ones, without
# touching its newly-set value.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats

# Data cells: a constant "1" flag for every existing data row.
$ws.Range("H2").Value = 1
$ws.Range("H3").Value = 1
$ws.Range("H4").Value = 1
